$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.527.78'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '2.483.46'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.45'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.62'
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('D13').Value = '2.867.28'
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.41'
$ws.Range('E15').Value = '  +7.34%  '
$ws.Range('D16').Value = '2.469.11'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.754'
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('D18').Value = '41.664.92'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.30'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.51'
$ws.Range('E21').Value = '  +5.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.15'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.24'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.80'
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.21'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.54'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.40'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.16'
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0754'
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.64'
$ws.Range('E43').Value = '  -6.84%  '
$ws.Range('D44').Value = '1.956.15'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0284'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.86'
$ws.Range('E47').Value = '  +2.06%  '
$ws.Range('D48').Value = '2.725.71'
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.07'
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.09'
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.13'
$ws.Range('E51').Value = '  -3.64%  '
